$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.040714492914588667
$ws.Range("B1").Value = 0.040714491002168091

$ws.Range("A2").Value = 0.060572896919914104
$ws.Range("B2").Value = -0.060572898682522539

$ws.Range("A3").Value = 0.021578416827811388
$ws.Range("B3").Value = -0.021578418624398493

$ws.Range("A4").Value = 0.033738693456479646
$ws.Range("B4").Value = -0.033738695254919397

$ws.Range("A5").Value = -0.0066223264963496242
$ws.Range("B5").Value = 0.0066223244945752789
